# Generate Report for Handoff
# Swaps the in-flight handoff GUID/hash for a fresh one and refreshes the
# associated timestamps, clearing the stale "Latest Target File" columns
# that no longer apply to the new handoff.

$wb = $excel.ActiveWorkbook

$oldGuid = "91171a96-5f7d-4641-8576-44b0058bbe25"
$newGuid = "e3607e7e-7a60-4fa4-8ce0-94bb847f5a1f"
$oldHash = "9800bf1d99ddd4eb7800e02a0ef556eabf3e0ef6"
$newHash = "f1015d7c5127d4729e7e7662dc936d340f028302"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-29 02:59:31"

$hOverview = $wsOverview.Hyperlinks
$hOverview.Delete()
$hOverview.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1684b09a54dcde7f52c067f3dd5e276c8aa47a75/e2e/$newGuid.md", "", "", "e2e\$newGuid.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-29 02:59:26"
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$hZh = $wsZh.Hyperlinks
$hZh.Delete()
$hZh.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1684b09a54dcde7f52c067f3dd5e276c8aa47a75/e2e/$newGuid.md", "", "", "$newGuid.md")

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-29 02:59:31"
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$hDe = $wsDe.Hyperlinks
$hDe.Delete()
$hDe.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1684b09a54dcde7f52c067f3dd5e276c8aa47a75/e2e/$newGuid.md", "", "", "$newGuid.md")

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
